# Apply cryptos list update (prices + volume%) per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.445.61'
$ws.Range('E2').Value = '  -3.48%  '
$ws.Range('D3').Value = '1.995.33'
$ws.Range('E3').Value = '  -5.99%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.009'
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '329.54'
$ws.Range('E5').Value = '  -5.10%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.007'
$ws.Range('E6').Value = '  +0.07%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5009'
$ws.Range('E7').Value = '  -4.36%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4228'
$ws.Range('E8').Value = '  -5.58%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '53.39'
$ws.Range('E9').Value = '  -1.91%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.08931'
$ws.Range('E10').Value = '  -5.04%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.122'
$ws.Range('E11').Value = '  -5.04%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '23.25'
$ws.Range('E12').Value = '  -7.85%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '8.086'
$ws.Range('E13').Value = '  -7.09%  '
$ws.Range('D14').Value = '1.983.11'
$ws.Range('E14').Value = '  -5.58%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.522'
$ws.Range('E15').Value = '  -6.38%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '95.77'
$ws.Range('E16').Value = '  -6.42%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.008'
$ws.Range('E17').Value = '  -0.02%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001109'
$ws.Range('E18').Value = '  -5.01%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06617'
$ws.Range('E19').Value = '  -1.63%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '19.67'
$ws.Range('E20').Value = '  -8.33%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.007'
$ws.Range('E21').Value = '  +0.03%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.971'
$ws.Range('E22').Value = '  -5.94%  '
$ws.Range('D23').Value = '29.465.61'
$ws.Range('E23').Value = '  -3.48%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.89'
$ws.Range('E24').Value = '  -6.83%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.254'
$ws.Range('E25').Value = '  -3.14%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '158.25'
$ws.Range('E26').Value = '  -2.58%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.66'
$ws.Range('E27').Value = '  -7.27%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.556'
$ws.Range('E28').Value = '  -4.86%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.327'
$ws.Range('E29').Value = '  -8.87%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '127.82'
$ws.Range('E30').Value = '  -4.92%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.048'
$ws.Range('E31').Value = '  -9.96%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09944'
$ws.Range('E32').Value = '  -6.39%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.569'
$ws.Range('E33').Value = '  -11.66%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.847'
$ws.Range('E34').Value = '  -7.18%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.784'
$ws.Range('E35').Value = '  -4.54%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '9.572'
$ws.Range('E36').Value = '  -9.83%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02462'
$ws.Range('E37').Value = '  -7.36%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.06333'
$ws.Range('E38').Value = '  -7.76%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.286'
$ws.Range('E39').Value = '  -3.88%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.6522'
$ws.Range('E40').Value = '  -8.51%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '11.70'
$ws.Range('E41').Value = '  -7.72%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.2066'
$ws.Range('E42').Value = '  -8.00%  '
$ws.Range('E43').Value = '  +0.05%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.6330'
$ws.Range('E44').Value = '  -8.77%  '
$ws.Range('B45').Value = 'NEARProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.205'
$ws.Range('E45').Value = '  -7.93%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '13.38'
$ws.Range('E46').Value = '  -8.97%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.267'
$ws.Range('E47').Value = '  -4.52%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.528'
$ws.Range('E48').Value = '  -3.57%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.00000000330'
$ws.Range('E49').Value = '  -4.55%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06995'
$ws.Range('E50').Value = '  -3.37%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.137'
$ws.Range('E51').Value = '  -4.94%  '
